$d = $word.ActiveDocument

function New-WordXmlFragment($innerBodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $innerBodyXml + '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

# Step 1: locate the trailing " uppalapati" text and split its run into a
# plain " " run plus a separate "uppalapati" run wrapped in spell-check
# proofErr markers (as Word's background spell checker would do while typing).
$r1 = $d.Content.Duplicate
$r1.Find.Execute(" uppalapati", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$frag1 = '<w:p><w:r w:rsidR="009066E0"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>uppalapati</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$r1.InsertXML((New-WordXmlFragment $frag1))

# Step 2: add a new paragraph containing "Hii" (flagged as a spelling error)
# right after the paragraph that now ends in "uppalapati"
$lastParaEnd = $d.Paragraphs($d.Paragraphs.Count).Range.End
$r2 = $d.Range($lastParaEnd - 1, $lastParaEnd - 1)
$frag2 = '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
         '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Hii</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$r2.InsertXML((New-WordXmlFragment $frag2))

# Step 3: add a final, empty paragraph after the "Hii" paragraph
$lastParaEnd = $d.Paragraphs($d.Paragraphs.Count).Range.End
$r3 = $d.Range($lastParaEnd - 1, $lastParaEnd - 1)
$frag3 = '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>'
$r3.InsertXML((New-WordXmlFragment $frag3))

Write-Host "Final content:" $d.Content.Text
Write-Host "Paragraph count:" $d.Paragraphs.Count
